$d = $word.ActiveDocument

$d.Content.Find.Execute("Java, C, Python, Bash, Unix Shell, Ruby, Lisp, Prolog", $true, $false, $false, $false, $false, $true, 1, $false, "Java, C, Python, Bash/Unix Shell, Ruby, Lisp, Prolog", 2)

$d.Content.Find.Execute("Familiar with SQL, MySQL", $true, $false, $false, $false, $false, $true, 1, $false, "Familiar with MySQL and relational databases", 2)
